$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New header / help-text content for columns AA:AD (27-30), rows 4, 6 and 7.
# Row 5 ("Required or Optional" row) reuses the existing "Conditional" string.
# ---------------------------------------------------------------------------

# Row 4 - Field ID
$ws.Range("AA4").Value = "Personnel_Exp_Exceeding_Estimate__c"
$ws.Range("AB4").Value = "Personnel_Obligations_Pursuant_Estimate__c"
$ws.Range("AC4").Value = "Contract_Expenditures_Exceeding_Estimate__c"
$ws.Range("AD4").Value = "Contract_Obligations_Pursuant_Estimate__c"

# Row 5 - Required or Optional -> Conditional (matches existing cells, e.g. Z5)
$ws.Range("AA5").Value = "Conditional"
$ws.Range("AB5").Value = "Conditional"
$ws.Range("AC5").Value = "Conditional"
$ws.Range("AD5").Value = "Conditional"
# AA5 picks up a bottom-less border (matches the style used by e.g. Z6/AC6/AD6)
$ws.Range("AA5").Borders.Item(9).LineStyle = -4142

# Row 6 - Field Name
$ws.Range("AA6").Value = "Personnel Expenditures Exceeding Estimate"
$ws.Range("AB6").Value = "Personnel Obligations Pursuant to Estimate"
$ws.Range("AC6").Value = "Contract Expenditures Exceeding Estimate"
$ws.Range("AD6").Value = "Contract Obligations Pursuant to Estimate"
# AA6 reverts to the plain column-default look (no border/override) - pull that
# format from a far-away untouched cell in the same column.
$ws.Range("AA979").Copy() | Out-Null
$ws.Range("AA6").PasteSpecial(-4122) | Out-Null

# Row 7 - Help Text
$ws.Range("AA7").Value = "Total Expenditures Exceeding Estimate" + [char]10 + "Required If Subaward Type is Direct Payment "
$ws.Range("AB7").Value = "Total Reported Obligations Pursuant to Estimate" + [char]10 + "Required If Subaward Type is Direct Payment "
$ws.Range("AC7").Value = "Total Expenditures Exceeding Estimate" + [char]10 + "Required If Subaward Type is one of the following: " + [char]10 + [char]8220 + "Contract: Purchase Order" + [char]8221 + [char]10 + [char]8220 + "Contract: Delivery Order" + [char]8221 + [char]10 + [char]8220 + "Contract: Blanket Purchase Agreement" + [char]8221 + [char]10 + [char]8220 + "Contract: Definitive Contract" + [char]8221
$ws.Range("AD7").Value = "Total Reported Obligations Pursuant to Estimate" + [char]10 + "Required If Subaward Type is one of the following: " + [char]10 + [char]8220 + "Contract: Purchase Order" + [char]8221 + [char]10 + [char]8220 + "Contract: Delivery Order" + [char]8221 + [char]10 + [char]8220 + "Contract: Blanket Purchase Agreement" + [char]8221 + [char]10 + [char]8220 + "Contract: Definitive Contract" + [char]8221

# ---------------------------------------------------------------------------
# Column widths widen for AB:AD (28-30) to fit the new text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(28).ColumnWidth = 40.857142857142854
$ws.Columns.Item(29).ColumnWidth = 37.857142857142854
$ws.Columns.Item(30).ColumnWidth = 36.42857142857143

# ---------------------------------------------------------------------------
# Move the active selection to AD7.
# ---------------------------------------------------------------------------
$ws.Range("AD7").Select() | Out-Null
